# Adds a "Hygieia github link" row to the Definitions sheet, pointing at
# the project's GitHub repository, and leaves the workbook focused on
# that new row.
$wb = $excel.ActiveWorkbook

$wsDefs = $wb.Worksheets.Item("Definitions")

# --- Definitions sheet: fill in row 3 (label + URL) ---
# Shared-string table order mirrors the order cells are first populated,
# so write the label before the hyperlink's URL text.
$wsDefs.Range("B3").Value = "Hygieia github link"
$wsDefs.Range("B3").Style = "Normal"
$wsDefs.Range("B3").WrapText = $true

$wsDefs.Range("A3").Value = "https://github.com/capitalone/Hygieia"
$wsDefs.Hyperlinks.Add($wsDefs.Range("A3"), "https://github.com/capitalone/Hygieia") | Out-Null
$wsDefs.Range("A3").Style = "Hyperlink"

# The old, never-filled-in formatting placeholder in B4 goes away now
# that row 3 carries real content.
$wsDefs.Range("B4").Clear()

# --- Selection / active sheet / view state ---
$wsDefs.Activate()
$wsDefs.Range("B3").Select() | Out-Null

$wb.Save()
